$wb = $excel.ActiveWorkbook

# ===== Sheet "LP1912" (sheet1) =====
$ws1 = $wb.Worksheets.Item("LP1912")

# Header updates
$ws1.Cells.Item(2,1).Value = "Última actualización: 09:25:30"
$ws1.Cells.Item(3,1).Value = "Total filas: 143"

# Swap stops for rows 59/60 (column C)
$ws1.Cells.Item(59,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(60,3).Value = "16_SANTA ANA"

# Rewrite data rows 111-148 with updated schedule values
$ws1.Cells.Item(111,1).Value = "08:47:19"
$ws1.Cells.Item(111,2).Value = "09:23"
$ws1.Cells.Item(111,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(111,4).Value = 36
$ws1.Cells.Item(111,5).Value = "LP1912"
$ws1.Cells.Item(112,1).Value = "07:44:08"
$ws1.Cells.Item(112,2).Value = "09:23"
$ws1.Cells.Item(112,3).Value = "17_ROMERO"
$ws1.Cells.Item(112,4).Value = 99
$ws1.Cells.Item(112,5).Value = "LP1912"
$ws1.Cells.Item(113,1).Value = "07:57:27"
$ws1.Cells.Item(113,2).Value = "09:23"
$ws1.Cells.Item(113,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(113,4).Value = 86
$ws1.Cells.Item(113,5).Value = "LP1912"
$ws1.Cells.Item(114,1).Value = "07:44:08"
$ws1.Cells.Item(114,2).Value = "09:24"
$ws1.Cells.Item(114,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(114,4).Value = 100
$ws1.Cells.Item(114,5).Value = "LP1912"
$ws1.Cells.Item(115,1).Value = "08:16:48"
$ws1.Cells.Item(115,2).Value = "09:29"
$ws1.Cells.Item(115,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(115,4).Value = 73
$ws1.Cells.Item(115,5).Value = "LP1912"
$ws1.Cells.Item(116,1).Value = "07:44:08"
$ws1.Cells.Item(116,2).Value = "09:32"
$ws1.Cells.Item(116,3).Value = "15_ABASTO"
$ws1.Cells.Item(116,4).Value = 108
$ws1.Cells.Item(116,5).Value = "LP1912"
$ws1.Cells.Item(117,1).Value = "07:44:08"
$ws1.Cells.Item(117,2).Value = "09:33"
$ws1.Cells.Item(117,3).Value = "10_OLMOS"
$ws1.Cells.Item(117,4).Value = 109
$ws1.Cells.Item(117,5).Value = "LP1912"
$ws1.Cells.Item(118,1).Value = "08:33:47"
$ws1.Cells.Item(118,2).Value = "09:34"
$ws1.Cells.Item(118,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(118,4).Value = 61
$ws1.Cells.Item(118,5).Value = "LP1912"
$ws1.Cells.Item(119,1).Value = "08:54:42"
$ws1.Cells.Item(119,2).Value = "09:34"
$ws1.Cells.Item(119,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(119,4).Value = 40
$ws1.Cells.Item(119,5).Value = "LP1912"
$ws1.Cells.Item(120,1).Value = "08:47:19"
$ws1.Cells.Item(120,2).Value = "09:35"
$ws1.Cells.Item(120,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(120,4).Value = 48
$ws1.Cells.Item(120,5).Value = "LP1912"
$ws1.Cells.Item(121,1).Value = "08:47:19"
$ws1.Cells.Item(121,2).Value = "09:35"
$ws1.Cells.Item(121,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(121,4).Value = 48
$ws1.Cells.Item(121,5).Value = "LP1912"
$ws1.Cells.Item(122,1).Value = "07:44:08"
$ws1.Cells.Item(122,2).Value = "09:36"
$ws1.Cells.Item(122,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(122,4).Value = 112
$ws1.Cells.Item(122,5).Value = "LP1912"
$ws1.Cells.Item(123,1).Value = "08:16:48"
$ws1.Cells.Item(123,2).Value = "09:37"
$ws1.Cells.Item(123,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(123,4).Value = 81
$ws1.Cells.Item(123,5).Value = "LP1912"
$ws1.Cells.Item(124,1).Value = "08:16:48"
$ws1.Cells.Item(124,2).Value = "09:41"
$ws1.Cells.Item(124,3).Value = "215C_EL PATO"
$ws1.Cells.Item(124,4).Value = 85
$ws1.Cells.Item(124,5).Value = "LP1912"
$ws1.Cells.Item(125,1).Value = "08:33:47"
$ws1.Cells.Item(125,2).Value = "09:41"
$ws1.Cells.Item(125,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(125,4).Value = 68
$ws1.Cells.Item(125,5).Value = "LP1912"
$ws1.Cells.Item(126,1).Value = "07:44:08"
$ws1.Cells.Item(126,2).Value = "09:42"
$ws1.Cells.Item(126,3).Value = "215C_EL PATO"
$ws1.Cells.Item(126,4).Value = 118
$ws1.Cells.Item(126,5).Value = "LP1912"
$ws1.Cells.Item(127,1).Value = "07:57:27"
$ws1.Cells.Item(127,2).Value = "09:43"
$ws1.Cells.Item(127,3).Value = "14_ABASTO"
$ws1.Cells.Item(127,4).Value = 106
$ws1.Cells.Item(127,5).Value = "LP1912"
$ws1.Cells.Item(128,1).Value = "09:25:30"
$ws1.Cells.Item(128,2).Value = "09:46"
$ws1.Cells.Item(128,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(128,4).Value = 21
$ws1.Cells.Item(128,5).Value = "LP1912"
$ws1.Cells.Item(129,1).Value = "08:54:42"
$ws1.Cells.Item(129,2).Value = "09:52"
$ws1.Cells.Item(129,3).Value = "15_ABASTO"
$ws1.Cells.Item(129,4).Value = 58
$ws1.Cells.Item(129,5).Value = "LP1912"
$ws1.Cells.Item(130,1).Value = "08:54:42"
$ws1.Cells.Item(130,2).Value = "09:53"
$ws1.Cells.Item(130,3).Value = "10_OLMOS"
$ws1.Cells.Item(130,4).Value = 59
$ws1.Cells.Item(130,5).Value = "LP1912"
$ws1.Cells.Item(131,1).Value = "09:25:30"
$ws1.Cells.Item(131,2).Value = "09:58"
$ws1.Cells.Item(131,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(131,4).Value = 33
$ws1.Cells.Item(131,5).Value = "LP1912"
$ws1.Cells.Item(132,1).Value = "09:25:30"
$ws1.Cells.Item(132,2).Value = "10:03"
$ws1.Cells.Item(132,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(132,4).Value = 38
$ws1.Cells.Item(132,5).Value = "LP1912"
$ws1.Cells.Item(133,1).Value = "09:25:30"
$ws1.Cells.Item(133,2).Value = "10:04"
$ws1.Cells.Item(133,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(133,4).Value = 39
$ws1.Cells.Item(133,5).Value = "LP1912"
$ws1.Cells.Item(134,1).Value = "08:16:48"
$ws1.Cells.Item(134,2).Value = "10:10"
$ws1.Cells.Item(134,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(134,4).Value = 114
$ws1.Cells.Item(134,5).Value = "LP1912"
$ws1.Cells.Item(135,1).Value = "08:16:48"
$ws1.Cells.Item(135,2).Value = "10:12"
$ws1.Cells.Item(135,3).Value = "15_ABASTO"
$ws1.Cells.Item(135,4).Value = 116
$ws1.Cells.Item(135,5).Value = "LP1912"
$ws1.Cells.Item(136,1).Value = "09:25:30"
$ws1.Cells.Item(136,2).Value = "10:13"
$ws1.Cells.Item(136,3).Value = "10_OLMOS"
$ws1.Cells.Item(136,4).Value = 48
$ws1.Cells.Item(136,5).Value = "LP1912"
$ws1.Cells.Item(137,1).Value = "08:33:47"
$ws1.Cells.Item(137,2).Value = "10:21"
$ws1.Cells.Item(137,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(137,4).Value = 108
$ws1.Cells.Item(137,5).Value = "LP1912"
$ws1.Cells.Item(138,1).Value = "08:33:47"
$ws1.Cells.Item(138,2).Value = "10:22"
$ws1.Cells.Item(138,3).Value = "17_ROMERO"
$ws1.Cells.Item(138,4).Value = 109
$ws1.Cells.Item(138,5).Value = "LP1912"
$ws1.Cells.Item(139,1).Value = "09:25:30"
$ws1.Cells.Item(139,2).Value = "10:23"
$ws1.Cells.Item(139,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(139,4).Value = 58
$ws1.Cells.Item(139,5).Value = "LP1912"
$ws1.Cells.Item(140,1).Value = "08:33:47"
$ws1.Cells.Item(140,2).Value = "10:26"
$ws1.Cells.Item(140,3).Value = "215A_EL PATO"
$ws1.Cells.Item(140,4).Value = 113
$ws1.Cells.Item(140,5).Value = "LP1912"
$ws1.Cells.Item(141,1).Value = "08:54:42"
$ws1.Cells.Item(141,2).Value = "10:41"
$ws1.Cells.Item(141,3).Value = "17_ROMERO"
$ws1.Cells.Item(141,4).Value = 107
$ws1.Cells.Item(141,5).Value = "LP1912"
$ws1.Cells.Item(142,1).Value = "08:47:19"
$ws1.Cells.Item(142,2).Value = "10:42"
$ws1.Cells.Item(142,3).Value = "17_ROMERO"
$ws1.Cells.Item(142,4).Value = 115
$ws1.Cells.Item(142,5).Value = "LP1912"
$ws1.Cells.Item(143,1).Value = "08:47:19"
$ws1.Cells.Item(143,2).Value = "10:43"
$ws1.Cells.Item(143,3).Value = "14_ABASTO"
$ws1.Cells.Item(143,4).Value = 116
$ws1.Cells.Item(143,5).Value = "LP1912"
$ws1.Cells.Item(144,1).Value = "09:25:30"
$ws1.Cells.Item(144,2).Value = "10:53"
$ws1.Cells.Item(144,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(144,4).Value = 88
$ws1.Cells.Item(144,5).Value = "LP1912"
$ws1.Cells.Item(145,1).Value = "09:25:30"
$ws1.Cells.Item(145,2).Value = "11:02"
$ws1.Cells.Item(145,3).Value = "215C_EL PATO"
$ws1.Cells.Item(145,4).Value = 97
$ws1.Cells.Item(145,5).Value = "LP1912"
$ws1.Cells.Item(146,1).Value = "09:25:30"
$ws1.Cells.Item(146,2).Value = "11:06"
$ws1.Cells.Item(146,3).Value = "16_P MOR-167 Y 521"
$ws1.Cells.Item(146,4).Value = 101
$ws1.Cells.Item(146,5).Value = "LP1912"
$ws1.Cells.Item(147,1).Value = "09:25:30"
$ws1.Cells.Item(147,2).Value = "11:19"
$ws1.Cells.Item(147,3).Value = "86_EST CHICA-ESC AGRARIA"
$ws1.Cells.Item(147,4).Value = 114
$ws1.Cells.Item(147,5).Value = "LP1912"
$ws1.Cells.Item(148,1).Value = "09:25:30"
$ws1.Cells.Item(148,2).Value = "11:21"
$ws1.Cells.Item(148,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(148,4).Value = 116
$ws1.Cells.Item(148,5).Value = "LP1912"

# ===== Sheet "LP1912-215" (sheet2) =====
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2,1).Value = "Última actualización: 09:25:30"
$ws2.Cells.Item(3,1).Value = "Total filas: 20"

# New data row 25
$ws2.Cells.Item(25,1).Value = "09:25:30"
$ws2.Cells.Item(25,2).Value = "11:02"
$ws2.Cells.Item(25,3).Value = "215C_EL PATO"
$ws2.Cells.Item(25,4).Value = 97
$ws2.Cells.Item(25,5).Value = "LP1912"

# ===== Sheet "6203-6173" (sheet3) =====
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2,1).Value = "Última actualización: 09:25:30"
$ws3.Cells.Item(3,1).Value = "Total filas: 27"

# New data rows 31-32
$ws3.Cells.Item(31,1).Value = "09:25:30"
$ws3.Cells.Item(31,2).Value = "10:54"
$ws3.Cells.Item(31,3).Value = "215A_LA PLATA"
$ws3.Cells.Item(31,4).Value = 89
$ws3.Cells.Item(31,5).Value = "L6173"
$ws3.Cells.Item(32,1).Value = "09:25:30"
$ws3.Cells.Item(32,2).Value = "11:14"
$ws3.Cells.Item(32,3).Value = "215C_LA PLATA"
$ws3.Cells.Item(32,4).Value = 109
$ws3.Cells.Item(32,5).Value = "L6203"